$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Add the new log entries (rows 719-739, 743-748; rows 740-742 are intentionally left blank,
# matching gaps already present elsewhere in the log).

$ws.Cells.Item(719, 1).Value = 'Demo'
$ws.Cells.Item(719, 2).Value = 42661
$ws.Cells.Item(719, 3).Value = '1730'
$ws.Cells.Item(719, 4).Value = 'SSB'
$ws.Cells.Item(719, 5).Value = 'S126'

$ws.Cells.Item(720, 1).Value = 'Operator'
$ws.Cells.Item(720, 2).Value = 42661
$ws.Cells.Item(720, 3).Value = '1830'
$ws.Cells.Item(720, 4).Value = 'OSG'
$ws.Cells.Item(720, 5).Value = '1014G'
$ws.Cells.Item(720, 6).Value = 'Please remain on site and oversee rooms 1003, 1014H, 2008, 2009, 2010 recordings'
$ws.Rows.Item(720).RowHeight = 30

$ws.Cells.Item(721, 1).Value = 'Setup Skype Kit'
$ws.Cells.Item(721, 2).Value = 42661
$ws.Cells.Item(721, 3).Value = '1630'
$ws.Cells.Item(721, 4).Value = 'OSG'
$ws.Cells.Item(721, 5).Value = '1014G'
$ws.Cells.Item(721, 6).Value = 'Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L'
$ws.Rows.Item(721).RowHeight = 30

$ws.Cells.Item(722, 1).Value = 'Setup Skype Kit'
$ws.Cells.Item(722, 2).Value = 42661
$ws.Cells.Item(722, 3).Value = '1630'
$ws.Cells.Item(722, 4).Value = 'OSG'
$ws.Cells.Item(722, 5).Value = '1014J'
$ws.Cells.Item(722, 6).Value = 'Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L'
$ws.Rows.Item(722).RowHeight = 30

$ws.Cells.Item(723, 1).Value = 'Setup Skype Kit'
$ws.Cells.Item(723, 2).Value = 42661
$ws.Cells.Item(723, 3).Value = '1630'
$ws.Cells.Item(723, 4).Value = 'OSG'
$ws.Cells.Item(723, 5).Value = '1014K'
$ws.Cells.Item(723, 6).Value = 'Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L'
$ws.Rows.Item(723).RowHeight = 30

$ws.Cells.Item(724, 1).Value = 'Setup Skype Kit'
$ws.Cells.Item(724, 2).Value = 42661
$ws.Cells.Item(724, 3).Value = '1630'
$ws.Cells.Item(724, 4).Value = 'OSG'
$ws.Cells.Item(724, 5).Value = '2001'
$ws.Cells.Item(724, 6).Value = 'Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L'
$ws.Rows.Item(724).RowHeight = 30

$ws.Cells.Item(725, 1).Value = 'Setup Skype Kit'
$ws.Cells.Item(725, 2).Value = 42661
$ws.Cells.Item(725, 3).Value = '1630'
$ws.Cells.Item(725, 4).Value = 'OSG'
$ws.Cells.Item(725, 5).Value = '2002'
$ws.Cells.Item(725, 6).Value = 'Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L'
$ws.Rows.Item(725).RowHeight = 30

$ws.Cells.Item(726, 1).Value = 'Demo'
$ws.Cells.Item(726, 2).Value = 42661
$ws.Cells.Item(726, 3).Value = '1815'
$ws.Cells.Item(726, 4).Value = 'OSG'
$ws.Cells.Item(726, 5).Value = '1014G'
$ws.Cells.Item(726, 6).Value = 'Make sure recording is going well. Demo AV equipment. '

$ws.Cells.Item(727, 1).Value = 'Demo'
$ws.Cells.Item(727, 2).Value = 42661
$ws.Cells.Item(727, 3).Value = '1815'
$ws.Cells.Item(727, 4).Value = 'OSG'
$ws.Cells.Item(727, 5).Value = '1014J'
$ws.Cells.Item(727, 6).Value = 'Make sure recording is going well. Demo AV equipment. '

$ws.Cells.Item(728, 1).Value = 'Demo'
$ws.Cells.Item(728, 2).Value = 42661
$ws.Cells.Item(728, 3).Value = '1815'
$ws.Cells.Item(728, 4).Value = 'OSG'
$ws.Cells.Item(728, 5).Value = '1014K'
$ws.Cells.Item(728, 6).Value = 'Make sure recording is going well. Demo AV equipment. '

$ws.Cells.Item(729, 1).Value = 'Demo'
$ws.Cells.Item(729, 2).Value = 42661
$ws.Cells.Item(729, 3).Value = '1815'
$ws.Cells.Item(729, 4).Value = 'OSG'
$ws.Cells.Item(729, 5).Value = '2001'
$ws.Cells.Item(729, 6).Value = 'Make sure recording is going well. Demo AV equipment. '

$ws.Cells.Item(730, 1).Value = 'Demo'
$ws.Cells.Item(730, 2).Value = 42661
$ws.Cells.Item(730, 3).Value = '1815'
$ws.Cells.Item(730, 4).Value = 'OSG'
$ws.Cells.Item(730, 5).Value = '2002'
$ws.Cells.Item(730, 6).Value = 'Make sure recording is going well. Demo AV equipment. '

$ws.Cells.Item(731, 1).Value = 'Pickup Skype Kit'
$ws.Cells.Item(731, 2).Value = 42661
$ws.Cells.Item(731, 3).Value = '2150'
$ws.Cells.Item(731, 4).Value = 'OSG'
$ws.Cells.Item(731, 5).Value = '1014G'
$ws.Cells.Item(731, 6).Value = 'Return web cam and tripod to OSG 1014L'

$ws.Cells.Item(732, 1).Value = 'Pickup Skype Kit'
$ws.Cells.Item(732, 2).Value = 42661
$ws.Cells.Item(732, 3).Value = '2150'
$ws.Cells.Item(732, 4).Value = 'OSG'
$ws.Cells.Item(732, 5).Value = '1014J'
$ws.Cells.Item(732, 6).Value = 'Return web cam and tripod to OSG 1014L'

$ws.Cells.Item(733, 1).Value = 'Pickup Skype Kit'
$ws.Cells.Item(733, 2).Value = 42661
$ws.Cells.Item(733, 3).Value = '2150'
$ws.Cells.Item(733, 4).Value = 'OSG'
$ws.Cells.Item(733, 5).Value = '1014K'
$ws.Cells.Item(733, 6).Value = 'Return web cam and tripod to OSG 1014L'

$ws.Cells.Item(734, 1).Value = 'Pickup Skype Kit'
$ws.Cells.Item(734, 2).Value = 42661
$ws.Cells.Item(734, 3).Value = '2150'
$ws.Cells.Item(734, 4).Value = 'OSG'
$ws.Cells.Item(734, 5).Value = '2001'
$ws.Cells.Item(734, 6).Value = 'Return web cam and tripod to OSG 1014L'

$ws.Cells.Item(735, 1).Value = 'Pickup Skype Kit'
$ws.Cells.Item(735, 2).Value = 42661
$ws.Cells.Item(735, 3).Value = '2150'
$ws.Cells.Item(735, 4).Value = 'OSG'
$ws.Cells.Item(735, 5).Value = '2002'
$ws.Cells.Item(735, 6).Value = 'Return web cam and tripod to OSG 1014L'

$ws.Cells.Item(736, 1).Value = 'SCLD Student Event'
$ws.Cells.Item(736, 2).Value = 42661
$ws.Cells.Item(736, 3).Value = '1830'
$ws.Cells.Item(736, 4).Value = 'FC'
$ws.Cells.Item(736, 5).Value = '104'
$ws.Cells.Item(736, 6).Value = 'INC000000733528'

$ws.Cells.Item(737, 1).Value = 'SCLD Student Logout'
$ws.Cells.Item(737, 2).Value = 42661
$ws.Cells.Item(737, 3).Value = '2030'
$ws.Cells.Item(737, 4).Value = 'FC'
$ws.Cells.Item(737, 5).Value = '104'
$ws.Cells.Item(737, 6).Value = 'INC000000733528'

$ws.Cells.Item(738, 1).Value = 'SCLD Student Event'
$ws.Cells.Item(738, 2).Value = 42661
$ws.Cells.Item(738, 3).Value = '1900'
$ws.Cells.Item(738, 4).Value = 'ACE'
$ws.Cells.Item(738, 5).Value = '004'
$ws.Cells.Item(738, 6).Value = 'INC000000737081'

$ws.Cells.Item(739, 1).Value = 'SCLD Student Logout'
$ws.Cells.Item(739, 2).Value = 42661
$ws.Cells.Item(739, 3).Value = '2100'
$ws.Cells.Item(739, 4).Value = 'ACE'
$ws.Cells.Item(739, 5).Value = '004'
$ws.Cells.Item(739, 6).Value = 'INC000000737081'

$ws.Cells.Item(743, 1).Value = 'Pickup Skype Kit'
$ws.Cells.Item(743, 2).Value = 42662
$ws.Cells.Item(743, 3).Value = '1630'
$ws.Cells.Item(743, 4).Value = 'ACE'
$ws.Cells.Item(743, 5).Value = '004'
$ws.Cells.Item(743, 6).Value = 'Return to ACE 015'

$ws.Cells.Item(744, 1).Value = 'Pickup PC'
$ws.Cells.Item(744, 2).Value = 42662
$ws.Cells.Item(744, 3).Value = '1730'
$ws.Cells.Item(744, 4).Value = 'ATK'
$ws.Cells.Item(744, 5).Value = '005'
$ws.Cells.Item(744, 6).Value = 'Return to DB 0003 !!!!!'

$ws.Cells.Item(745, 1).Value = 'Pickup Projector'
$ws.Cells.Item(745, 2).Value = 42662
$ws.Cells.Item(745, 3).Value = '1730'
$ws.Cells.Item(745, 4).Value = 'ATK'
$ws.Cells.Item(745, 5).Value = '005'
$ws.Cells.Item(745, 6).Value = 'Return to ATK 003C'

$ws.Cells.Item(746, 1).Value = 'AV Shutdown'
$ws.Cells.Item(746, 2).Value = 42662
$ws.Cells.Item(746, 3).Value = '2100'
$ws.Cells.Item(746, 4).Value = 'SSB'
$ws.Cells.Item(746, 5).Value = 'W141'

$ws.Cells.Item(747, 1).Value = 'Setup PC'
$ws.Cells.Item(747, 2).Value = 42662
$ws.Cells.Item(747, 3).Value = '1630'
$ws.Cells.Item(747, 4).Value = 'HNE'
$ws.Cells.Item(747, 5).Value = '105'
$ws.Cells.Item(747, 6).Value = 'Equipment from HNES 003'

$ws.Cells.Item(748, 1).Value = 'Pickup PC'
$ws.Cells.Item(748, 2).Value = 42662
$ws.Cells.Item(748, 3).Value = '1730'
$ws.Cells.Item(748, 4).Value = 'HNE'
$ws.Cells.Item(748, 5).Value = '105'
$ws.Cells.Item(748, 6).Value = 'Return equipment to HNES 003'

# Restore the sheet view / selection to match where the log was last edited.
$ws.Activate()
$ws.Range("C752").Select()
